# "Change: allow WDPR (timber) on flat wagons"
#
# The cargo-classes sheet had an AutoFilter active that showed only the
# blank rows of the "Boxcar w Mail+Armoured" column (AL) -- that hid
# every cargo row that already had a wagon-class assignment. Clear that
# filter so every row is visible again, then mark WDPR (row 67) as
# allowed on "Flat/stake wagon" (column AF) by putting a 1 in AF67.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the AutoFilter criteria (removes the <filterColumn> predicate and
# un-hides every row that the filter had hidden).
[void]$ws.ShowAllData()

# Allow WDPR (Wood Products) on the flat/stake wagon class.
$ws.Range("AF67").Value = 1

# Leave the selection where the author's commit left it.
[void]$ws.Range("AF66").Select()
